# Daily update at 8 AM UTC
# Appends the next day's row of win counts to the "Wins Over Time" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row after the existing data (row 76 -> new row 77).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Copy the date cell's formatting (style) from the previous row down to the
# new row, then overwrite the values for the new day.
$ws.Range("A" + $lastRow).Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item($newRow, 1).Value = 46026
$ws.Cells.Item($newRow, 2).Value = 172
$ws.Cells.Item($newRow, 3).Value = 181
$ws.Cells.Item($newRow, 4).Value = 174
